$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 23.02.2022 05:30"

# Fix D5: was stored as text "+0.2", should be numeric 0.2
$ws.Range("D5").Value = 0.2

# Fix E5: was stored as text "2022-02-23 05:15:07", should be numeric date serial
# with the same date style (s="2") as the other E column cells (yyyy-mm-dd hh:mm:ss)
$ws.Range("E5").Value = 44615.21883101852
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
